# "Tabela pinos.xlsx" edit: update 74HC595 bit-assignment list, refresh the
# active-sheet/selection view state, and drop the two now-unused "Candeeiro"
# strings in favour of two new fixtures appended at the end of the list.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ESP32")
$ws2 = $wb.Worksheets.Item("74HC595")

# --- Update the 74HC595 output-bit table (column B, rows 2-9) -------------
# Net effect of the diff: "Candeeiro Sala" drops out (rows below shift up
# one) and two brand-new fixtures ("Armario Hall" / "Guarda-vestidos") are
# appended as the new last two rows; "Candeeiro Quarto" also drops out of
# the shared-string table entirely as a result.
$ws2.Range("B5").Value = "Luz Cozinha"
$ws2.Range("B6").Value = "Luz Quarto"
$ws2.Range("B7").Value = "Luz Casa de banho"
$ws2.Range("B8").Value = "Armário Hall"
$ws2.Range("B9").Value = "Guarda-vestidos"

# --- View state: ESP32 becomes the active/selected sheet ------------------
$ws2.Select()
$ws2.Range("J9").Select()

$ws1.Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C17").Select()
